# aggiornamento fino a 27/05
# Adds daily Covid-19 data rows 256-269 (2021-05-14 .. 2021-05-27)
# to Sheet1, mirroring the formatting of the last existing row (255).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 255
$firstNewRow = 256
$lastNewRow = 269

# Copy the formatting (styles/borders/number format) of the last existing
# row down into each of the new rows before writing values into them, so
# that column A keeps its date style.
$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)
for ($i = $firstNewRow; $i -le $lastNewRow; $i++) {
    $dstRow = $ws.Range("A" + $i + ":D" + $i)
    $srcRow.Copy($dstRow)
}

# r, date serial (A), nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$data = @(
    @(256, 44330, 0, 16, 91.06949741021117),
    @(257, 44331, 1, 10, 56.91843588138198),
    @(258, 44332, 0, 8, 45.53474870510559),
    @(259, 44333, 2, 9, 51.22659229324378),
    @(260, 44334, 1, 5, 28.45921794069099),
    @(261, 44335, 1, 6, 34.15106152882919),
    @(262, 44336, 4, 9, 51.22659229324378),
    @(263, 44337, 0, 9, 51.22659229324378),
    @(264, 44338, 1, 9, 51.22659229324378),
    @(265, 44339, 0, 9, 51.22659229324378),
    @(266, 44340, 1, 8, 45.53474870510559),
    @(267, 44341, 0, 7, 39.84290511696739),
    @(268, 44342, 2, 8, 45.53474870510559),
    @(269, 44343, 2, 6, 34.15106152882919)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
